$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: expiry-date column label format change ---
$ws.Range("I1").Value = "Expiry Date (dd-mm-yyyy)"

# --- Row 2 ---
$ws.Range("L2").Value = 35
$ws.Range("M2").Value = 3
$ws.Range("L2").NumberFormat = "General"
$ws.Range("M2").NumberFormat = "General"

# --- Row 3 ---
$ws.Range("L3").Value = 432
$ws.Range("M3").Value = 4
$ws.Range("L3").NumberFormat = "General"
$ws.Range("M3").NumberFormat = "General"

# --- Row 4 ---
$ws.Range("B4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = $null
$ws.Range("L4").Value = 35
$ws.Range("M4").Value = 3
$ws.Range("L4").NumberFormat = "General"
$ws.Range("M4").NumberFormat = "General"

# --- Row 5 ---
$ws.Range("E5").Value = $null
$ws.Range("F5").Value = $null
$ws.Range("L5").Value = $null
$ws.Range("M5").Value = 4
$ws.Range("L5").NumberFormat = "General"
$ws.Range("M5").NumberFormat = "General"

# --- Row 6 ---
$ws.Range("K6").Value = $null
$ws.Range("L6").Value = 35
$ws.Range("M6").Value = 3
$ws.Range("L6").NumberFormat = "General"
$ws.Range("M6").NumberFormat = "General"

# --- Row 7 ---
$ws.Range("C7").Value = $null
$ws.Range("I7").Value = $null
$ws.Range("L7").Value = 432
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = $null
$ws.Range("L7").NumberFormat = "General"
$ws.Range("M7").NumberFormat = "General"

# --- Selection state ---
$ws.Range("F18").Select()
